$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# E23: 480 -> 360
$ws.Range("E23").Value = 360

# Row 24: new time-log entry (12월 8일)
$ws.Range("A24").Value = "12월 8일"
$a24chars = $ws.Range("A24").Characters(4, 3)
$a24chars.Font.Name = "돋움"
$a24chars.Font.Size = 10
$ws.Range("B24").Value = 0.375
$ws.Range("C24").Value = 0.66666666666666663
$ws.Range("D24").Value = 60
$ws.Range("E24").Value = 420
$ws.Range("F24").Value = "Nodejs 코딩"

# Row 25: new time-log entry (12월 9일)
$ws.Range("A25").Value = "12월 9일"
$a25chars = $ws.Range("A25").Characters(4, 3)
$a25chars.Font.Name = "돋움"
$a25chars.Font.Size = 10
$ws.Range("B25").Value = 0.375
$ws.Range("C25").Value = 0.79166666666666663
$ws.Range("D25").Value = 100
$ws.Range("E25").Value = 600
$ws.Range("F25").Value = "Nodejs 코딩"

# Update the active selection to F24, matching the saved view state
$ws.Activate()
$ws.Range("F24").Select()
